$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.201.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.607.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "302.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3767"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3652"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.74"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.25%  "
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.271"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.15%  "
$ws.Range("E12").Value = "  -4.28%  "
$ws.Range("E13").Value = "  -3.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.574"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.571"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001268"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.608.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06780"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.05%  "
$ws.Range("E20").Value = "  -7.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.583"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.226.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.351"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.917"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.45%  "
$ws.Range("E27").Value = "  -4.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.254"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.418"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.968"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.788.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9794"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07739"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02786"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.289"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2555"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.92%  "
$ws.Range("E39").Value = "  -7.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08867"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.393"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7160"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6609"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.99%  "
$ws.Range("E46").Value = "  -6.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  -2.52%  "
$ws.Range("E49").Value = "  -3.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.172"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.46%  "
